$wb = $excel.ActiveWorkbook

# Update numeric cells in the Leve profit sheets per data refresh.
# Each block below updates the H/I/J/K/L/M/N cells for one leve row;
# cells not listed as changed are left untouched. Two rows (GSM!N141 and
# WVR!M39) have their trailing cell cleared entirely (no successor value).

$ws = $wb.Worksheets.Item("ALC")
# Row 138 (date 44169)
$ws.Range("H138").Value = 3495.8518
$ws.Range("I138").Value = 2946.44
$ws.Range("K138").Value = 8839.32
$ws.Range("M138").Value = -3699.32

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (date 43999)
$ws.Range("H61").Value = 2535.2432
$ws.Range("I61").Value = 1880.7407
$ws.Range("K61").Value = 1880.7407
$ws.Range("M61").Value = -1668.7407
# Row 110 (date 27708)
$ws.Range("H110").Value = 2094.3125
$ws.Range("I110").Value = 885.8461
$ws.Range("K110").Value = 885.8461
$ws.Range("M110").Value = 1159.1539
# Row 121 (date 26285)
$ws.Range("H121").Value = 49999
$ws.Range("J121").Value = 49999
$ws.Range("L121").Value = 49999
$ws.Range("N121").Value = -53493
# Row 132 (date 43997)
$ws.Range("H132").Value = 2593.318
$ws.Range("I132").Value = 2241.4905
$ws.Range("K132").Value = 6724.4715
$ws.Range("M132").Value = -4194.4715
# Row 136 (date 43999)
$ws.Range("H136").Value = 2535.2432
$ws.Range("I136").Value = 1880.7407
$ws.Range("K136").Value = 5642.2221
$ws.Range("M136").Value = -3092.2221

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (date 14149)
$ws.Range("H20").Value = 22734660
$ws.Range("I20").Value = 31258274
$ws.Range("K20").Value = 31258274
$ws.Range("M20").Value = -31258027
# Row 99 (date 19943)
$ws.Range("H99").Value = 43818.6
$ws.Range("I99").Value = 65236.375
$ws.Range("K99").Value = 65236.375
$ws.Range("M99").Value = -63738.375
# Row 107 (date 27706)
$ws.Range("H107").Value = 3345881.5
$ws.Range("I107").Value = 5129384.5
$ws.Range("K107").Value = 5129384.5
$ws.Range("M107").Value = -5127464.5
# Row 134 (date 43998)
$ws.Range("H134").Value = 2659.6445
$ws.Range("I134").Value = 2349.6052
$ws.Range("K134").Value = 7048.8156
$ws.Range("M134").Value = -4513.8156

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (date 44023)
$ws.Range("H31").Value = 3497.0393
$ws.Range("J31").Value = 6810.4546
$ws.Range("L31").Value = 6810.4546
$ws.Range("N31").Value = -7400.4546
# Row 34 (date 44023)
$ws.Range("H34").Value = 3497.0393
$ws.Range("J34").Value = 6810.4546
$ws.Range("L34").Value = 6810.4546
$ws.Range("N34").Value = -7214.4546
# Row 58 (date 44021)
$ws.Range("H58").Value = 2599.2
$ws.Range("I58").Value = 1373.5
$ws.Range("K58").Value = 1373.5
$ws.Range("M58").Value = -1170.5
# Row 122 (date 36196)
$ws.Range("H122").Value = 2498.9092
$ws.Range("I122").Value = 2843.4736
$ws.Range("J122").Value = 2031.2858
$ws.Range("K122").Value = 8530.4208
$ws.Range("L122").Value = 6093.857400000001
$ws.Range("M122").Value = -6080.4208
$ws.Range("N122").Value = -10993.8574
# Row 132 (date 44019)
$ws.Range("H132").Value = 26317564
$ws.Range("I132").Value = 33334382
$ws.Range("K132").Value = 100003146
$ws.Range("M132").Value = -100000616
# Row 136 (date 44021)
$ws.Range("H136").Value = 2599.2
$ws.Range("I136").Value = 1373.5
$ws.Range("K136").Value = 4120.5
$ws.Range("M136").Value = -1570.5

$ws = $wb.Worksheets.Item("CUL")
# Row 18 (date 36056)
$ws.Range("H18").Value = 1740
$ws.Range("I18").Value = 1740
$ws.Range("K18").Value = 5220
$ws.Range("M18").Value = -5051
# Row 56 (date 10146)
$ws.Range("H56").Value = 16085.556
$ws.Range("I56").Value = 16085.556
$ws.Range("K56").Value = 16085.556
$ws.Range("M56").Value = -15555.556
# Row 107 (date 27838)
$ws.Range("H107").Value = 252.28572
$ws.Range("J107").Value = 252.28572
$ws.Range("L107").Value = 756.85716
$ws.Range("N107").Value = -4596.85716
# Row 116 (date 27866)
$ws.Range("H116").Value = 115638.5
$ws.Range("I116").Value = 226492
$ws.Range("K116").Value = 679476
$ws.Range("M116").Value = -676034
# Row 117 (date 27870)
$ws.Range("H117").Value = 2065.818
$ws.Range("J117").Value = 1815.5
$ws.Range("L117").Value = 5446.5
$ws.Range("N117").Value = -12330.5
# Row 129 (date 36054)
$ws.Range("H129").Value = 1939.75
$ws.Range("I129").Value = 1939.75
$ws.Range("K129").Value = 5819.25
$ws.Range("M129").Value = -819.25

$ws = $wb.Worksheets.Item("GSM")
# Row 97 (date 19940)
$ws.Range("H97").Value = 3118.6
$ws.Range("I97").Value = 1398.25
$ws.Range("K97").Value = 1398.25
$ws.Range("M97").Value = -902.25
# Row 122 (date 36182)
$ws.Range("H122").Value = 3942.0908
$ws.Range("I122").Value = 3191.0667
$ws.Range("J122").Value = 5551.4287
$ws.Range("K122").Value = 9573.2001
$ws.Range("L122").Value = 16654.2861
$ws.Range("M122").Value = -7123.2001
$ws.Range("N122").Value = -21554.2861
# Row 126 (date 36184)
$ws.Range("H126").Value = 4375.1
$ws.Range("J126").Value = 12974
$ws.Range("L126").Value = 38922
$ws.Range("N126").Value = -43862
# Row 132 (date 44008)
$ws.Range("H132").Value = 2469.1
$ws.Range("I132").Value = 1946.5
$ws.Range("K132").Value = 5839.5
$ws.Range("M132").Value = -3309.5
# Row 141 (date 42504)
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (date 36249)
$ws.Range("H7").Value = 5999.909
$ws.Range("I7").Value = 7318.4
$ws.Range("K7").Value = 7318.4
$ws.Range("M7").Value = -7206.4
# Row 22 (date 5277)
$ws.Range("H22").Value = 556.3333
$ws.Range("I22").Value = 417
$ws.Range("J22").Value = 667.8
$ws.Range("K22").Value = 417
$ws.Range("L22").Value = 667.8
$ws.Range("M22").Value = -122
$ws.Range("N22").Value = -1257.8
# Row 27 (date 5277)
$ws.Range("H27").Value = 556.3333
$ws.Range("I27").Value = 417
$ws.Range("J27").Value = 667.8
$ws.Range("K27").Value = 417
$ws.Range("L27").Value = 667.8
$ws.Range("M27").Value = -310
$ws.Range("N27").Value = -881.8
# Row 31 (date 3043)
$ws.Range("H31").Value = 2178.5
$ws.Range("I31").Value = 1400
$ws.Range("J31").Value = 2957
$ws.Range("K31").Value = 1400
$ws.Range("L31").Value = 2957
$ws.Range("M31").Value = -1152
$ws.Range("N31").Value = -3453
# Row 40 (date 36248)
$ws.Range("H40").Value = 3960.0476
$ws.Range("I40").Value = 4713.154
$ws.Range("J40").Value = 2736.25
$ws.Range("K40").Value = 4713.154
$ws.Range("L40").Value = 2736.25
$ws.Range("M40").Value = -4577.154
$ws.Range("N40").Value = -3008.25
# Row 61 (date 27740)
$ws.Range("H61").Value = 28239.54
$ws.Range("I61").Value = 3710.2
$ws.Range("K61").Value = 3710.2
$ws.Range("M61").Value = -3508.2
# Row 100 (date 19995)
$ws.Range("H100").Value = 1515.8334
$ws.Range("I100").Value = 1324.5
$ws.Range("J100").Value = 1898.5
$ws.Range("K100").Value = 1324.5
$ws.Range("L100").Value = 1898.5
$ws.Range("M100").Value = -783.5
$ws.Range("N100").Value = -2980.5
# Row 113 (date 27740)
$ws.Range("H113").Value = 28239.54
$ws.Range("I113").Value = 3710.2
$ws.Range("K113").Value = 3710.2
$ws.Range("M113").Value = -1540.2
# Row 122 (date 36247)
$ws.Range("H122").Value = 3340.125
$ws.Range("I122").Value = 3087.1155
$ws.Range("K122").Value = 9261.3465
$ws.Range("M122").Value = -6811.3465
# Row 126 (date 36249)
$ws.Range("H126").Value = 5999.909
$ws.Range("I126").Value = 7318.4
$ws.Range("K126").Value = 21955.2
$ws.Range("M126").Value = -19485.2
# Row 132 (date 44058)
$ws.Range("H132").Value = 4234.8203
$ws.Range("I132").Value = 2148.4348
$ws.Range("J132").Value = 7234
$ws.Range("K132").Value = 6445.3044
$ws.Range("L132").Value = 21702
$ws.Range("M132").Value = -3915.3044
$ws.Range("N132").Value = -26762
# Row 136 (date 44060)
$ws.Range("H136").Value = 3685.8845
$ws.Range("I136").Value = 2862.3044
$ws.Range("K136").Value = 8586.913199999999
$ws.Range("M136").Value = -6036.913199999999

$ws = $wb.Worksheets.Item("WVR")
# Row 39 (date 3106)
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
# Row 107 (date 27746)
$ws.Range("H107").Value = 790.619
$ws.Range("I107").Value = 836.1177
$ws.Range("J107").Value = 597.25
$ws.Range("K107").Value = 2508.3531
$ws.Range("L107").Value = 1791.75
$ws.Range("M107").Value = -588.3531000000003
$ws.Range("N107").Value = -5631.75
# Row 119 (date 26289)
$ws.Range("H119").Value = 85415.336
$ws.Range("J119").Value = 85415.336
$ws.Range("L119").Value = 85415.336
$ws.Range("N119").Value = -95091.336
# Row 122 (date 36208)
$ws.Range("H122").Value = 19232366
$ws.Range("I122").Value = 1570.875
$ws.Range("J122").Value = 50001636
$ws.Range("K122").Value = 4712.625
$ws.Range("L122").Value = 150004908
$ws.Range("M122").Value = -2262.625
$ws.Range("N122").Value = -150009808
# Row 140 (date 42506)
$ws.Range("H140").Value = 91587.25
$ws.Range("J140").Value = 91587.25
$ws.Range("L140").Value = 91587.25
$ws.Range("N140").Value = -101947.25
# Row 141 (date 42505)
$ws.Range("H141").Value = 69999
$ws.Range("J141").Value = 69999
$ws.Range("L141").Value = 69999
$ws.Range("N141").Value = -80359
